{"js": "const pairs = [\n  [\"2024-08-24 Saturday\", \"2024-08-25 Sunday\"],\n  [\"125\u00f74=\", \"517\u00f76=\"],\n  [\"367\u00f78=\", \"766\u00f72=\"],\n  [\"226\u00f76=\", \"930\u00f76=\"],\n  [\"397\u00f73=\", \"631\u00f77=\"],\n  [\"364\u00f77=\", \"623\u00f73=\"],\n  [\"827\u00f76=\", \"217\u00f73=\"],\n  [\"695\u00f74=\", \"701\u00f79=\"],\n  [\"351\u00f77=\", \"825\u00f77=\"],\n  [\"310\u00f76=\", \"277\u00f74=\"],\n  [\"660\u00f72=\", \"475\u00f75=\"],\n  [\"950\u00f72=\", \"462\u00f73=\"],\n  [\"305\u00f72=\", \"741\u00f74=\"],\n  [\"792\u00f77=\", \"318\u00f79=\"],\n  [\"912\u00f75=\", \"325\u00f75=\"],\n  [\"382\u00f79=\", \"951\u00f74=\"],\n  [\"524\u00f75=\", \"204\u00f76=\"],\n  [\"588\u00f77=\", \"443\u00f73=\"],\n  [\"672\u00f78=\", \"152\u00f72=\"],\n  [\"993\u00f74=\", \"247\u00f76=\"],\n  [\"231\u00f73=\", \"176\u00f72=\"],\n  [\"914\u00f76=\", \"984\u00f72=\"],\n  [\"759\u00f73=\", \"618\u00f77=\"],\n  [\"742\u00f79=\", \"697\u00f78=\"],\n  [\"113\u00f78=\", \"286\u00f78=\"],\n  [\"539\u00f76=\", \"621\u00f75=\"],\n];\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-08-24 Saturday\", \"2024-08-25 Sunday\"),\n  @(\"125\u00f74=\", \"517\u00f76=\"),\n  @(\"367\u00f78=\", \"766\u00f72=\"),\n  @(\"226\u00f76=\", \"930\u00f76=\"),\n  @(\"397\u00f73=\", \"631\u00f77=\"),\n  @(\"364\u00f77=\", \"623\u00f73=\"),\n  @(\"827\u00f76=\", \"217\u00f73=\"),\n  @(\"695\u00f74=\", \"701\u00f79=\"),\n  @(\"351\u00f77=\", \"825\u00f77=\"),\n  @(\"310\u00f76=\", \"277\u00f74=\"),\n  @(\"660\u00f72=\", \"475\u00f75=\"),\n  @(\"950\u00f72=\", \"462\u00f73=\"),\n  @(\"305\u00f72=\", \"741\u00f74=\"),\n  @(\"792\u00f77=\", \"318\u00f79=\"),\n  @(\"912\u00f75=\", \"325\u00f75=\"),\n  @(\"382\u00f79=\", \"951\u00f74=\"),\n  @(\"524\u00f75=\", \"204\u00f76=\"),\n  @(\"588\u00f77=\", \"443\u00f73=\"),\n  @(\"672\u00f78=\", \"152\u00f72=\"),\n  @(\"993\u00f74=\", \"247\u00f76=\"),\n  @(\"231\u00f73=\", \"176\u00f72=\"),\n  @(\"914\u00f76=\", \"984\u00f72=\"),\n  @(\"759\u00f73=\", \"618\u00f77=\"),\n  @(\"742\u00f79=\", \"697\u00f78=\"),\n  @(\"113\u00f78=\", \"286\u00f78=\"),\n  @(\"539\u00f76=\", \"621\u00f75=\"),\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $ok = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $ok) {\n    throw \"Replace failed for: $oldText\"\n  }\n}\n\nWrite-Output \"done\""}
